# "Add files via upload" - re-upload of the CMIP6 workbook from a different
# machine/Excel build. The only substantive content change is in Sheet1:
# cell F2's label switches from "CN" to "CNP" (a new shared string "CNP" is
# introduced). The saved file also reflects the view state at save time
# (zoom level and the active cell selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2 from "CN" to "CNP" (adds a new shared string entry).
$ws.Range("F2").Value = "CNP"

# Reflect the view state captured in the saved workbook: zoomed to 137%
# with I3 as the active/selected cell.
$excel.ActiveWindow.Zoom = 137
$null = $ws.Range("I3").Select()
